$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Cell values (row 1 headers) ---
$ws.Range("A1").Value = "所屬一級單位"
$ws.Range("B1").Value = "所屬系所部門"
$ws.Range("C1").Value = "姓名"
$ws.Range("D1").Value = "身分 (學士、碩士或博士班）"
$ws.Range("E1").Value = "國籍"
$ws.Range("F1").Value = "開始時間"
$ws.Range("G1").Value = "結束時間"
$ws.Range("H1").Value = "備註"

# --- Black bold header style (A1:C1 and E1:H1), vertical-center aligned ---
$blackHeader1 = $ws.Range("A1:C1")
$blackHeader1.Font.Bold = $true
$blackHeader1.Font.Size = 12
$blackHeader1.Font.Name = "新細明體"
$blackHeader1.Font.Family = 1
$blackHeader1.VerticalAlignment = -4108

$blackHeader2 = $ws.Range("E1:H1")
$blackHeader2.Font.Bold = $true
$blackHeader2.Font.Size = 12
$blackHeader2.Font.Name = "新細明體"
$blackHeader2.Font.Family = 1
$blackHeader2.VerticalAlignment = -4108

# --- D1: bold blue header style, vertical-center aligned ---
$d1 = $ws.Range("D1")
$d1.Font.Bold = $true
$d1.Font.Size = 12
$d1.Font.Name = "新細明體"
$d1.Font.Family = 1
$d1.Font.Color = 16711680
$d1.VerticalAlignment = -4108

# --- Rich text: leading part of A1 / B1 shown in red (rest keeps the
#     black bold cell-level font) ---
$ws.Range("A1").Characters(1, 4).Font.Color = 255
$ws.Range("B1").Characters(1, 2).Font.Color = 255

# --- Column widths ---
$ws.Columns.Item(1).ColumnWidth = 16.142857142857142
$ws.Columns.Item(2).ColumnWidth = 13.142857142857142
$ws.Columns.Item(4).ColumnWidth = 30

# --- Selection moves to D7 ---
$ws.Range("D7").Activate()

Write-Host "done"
